$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null
$ws.Range("R4").Value = 2021
$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null
$ws.Range("R5").Value = 11.9
$ws.Range("Q6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null
$ws.Range("R6").Value = 13.1
$ws.Range("Q7").Copy() | Out-Null
$ws.Range("R7").PasteSpecial(-4122) | Out-Null
$ws.Range("R7").Value = 10.6
$ws.Range("Q8").Copy() | Out-Null
$ws.Range("R8").PasteSpecial(-4122) | Out-Null
$ws.Range("R8").Value = 11
$ws.Range("Q9").Copy() | Out-Null
$ws.Range("R9").PasteSpecial(-4122) | Out-Null
$ws.Range("R9").Value = 10
$ws.Range("Q10").Copy() | Out-Null
$ws.Range("R10").PasteSpecial(-4122) | Out-Null
$ws.Range("R10").Value = 12
$ws.Range("Q11").Copy() | Out-Null
$ws.Range("R11").PasteSpecial(-4122) | Out-Null
$ws.Range("R11").Value = 10.199999999999999
$ws.Range("Q12").Copy() | Out-Null
$ws.Range("R12").PasteSpecial(-4122) | Out-Null
$ws.Range("R12").Value = 10.5
$ws.Range("Q13").Copy() | Out-Null
$ws.Range("R13").PasteSpecial(-4122) | Out-Null
$ws.Range("R13").Value = 10
$ws.Range("Q14").Copy() | Out-Null
$ws.Range("R14").PasteSpecial(-4122) | Out-Null
$ws.Range("R14").Value = 19.399999999999999
$ws.Range("Q15").Copy() | Out-Null
$ws.Range("R15").PasteSpecial(-4122) | Out-Null
$ws.Range("R15").Value = 22.3
$ws.Range("Q16").Copy() | Out-Null
$ws.Range("R16").PasteSpecial(-4122) | Out-Null
$ws.Range("R16").Value = 16.399999999999999
$ws.Range("Q17").Copy() | Out-Null
$ws.Range("R17").PasteSpecial(-4122) | Out-Null
$ws.Range("R17").Value = 9.4
$ws.Range("Q18").Copy() | Out-Null
$ws.Range("R18").PasteSpecial(-4122) | Out-Null
$ws.Range("R18").Value = 11.4
$ws.Range("Q19").Copy() | Out-Null
$ws.Range("R19").PasteSpecial(-4122) | Out-Null
$ws.Range("R19").Value = 7.3
$ws.Range("Q20").Copy() | Out-Null
$ws.Range("R20").PasteSpecial(-4122) | Out-Null
$ws.Range("R20").Value = 3.1
$ws.Range("Q21").Copy() | Out-Null
$ws.Range("R21").PasteSpecial(-4122) | Out-Null
$ws.Range("R21").Value = 2.9
$ws.Range("Q22").Copy() | Out-Null
$ws.Range("R22").PasteSpecial(-4122) | Out-Null
$ws.Range("R22").Value = 3.4
$ws.Range("Q23").Copy() | Out-Null
$ws.Range("R23").PasteSpecial(-4122) | Out-Null
$ws.Range("R23").Value = 15
$ws.Range("Q24").Copy() | Out-Null
$ws.Range("R24").PasteSpecial(-4122) | Out-Null
$ws.Range("R24").Value = 17.3
$ws.Range("Q25").Copy() | Out-Null
$ws.Range("R25").PasteSpecial(-4122) | Out-Null
$ws.Range("R25").Value = 12.7
$ws.Range("Q26").Copy() | Out-Null
$ws.Range("R26").PasteSpecial(-4122) | Out-Null
$ws.Range("R26").Value = 7.9
$ws.Range("Q27").Copy() | Out-Null
$ws.Range("R27").PasteSpecial(-4122) | Out-Null
$ws.Range("R27").Value = 8.4
$ws.Range("Q28").Copy() | Out-Null
$ws.Range("R28").PasteSpecial(-4122) | Out-Null
$ws.Range("R28").Value = 7.4
$ws.Range("Q29").Copy() | Out-Null
$ws.Range("R29").PasteSpecial(-4122) | Out-Null
$ws.Range("R29").Value = 15.2
$ws.Range("Q30").Copy() | Out-Null
$ws.Range("R30").PasteSpecial(-4122) | Out-Null
$ws.Range("R30").Value = 17.600000000000001
$ws.Range("Q31").Copy() | Out-Null
$ws.Range("R31").PasteSpecial(-4122) | Out-Null
$ws.Range("R31").Value = 12.6
$ws.Range("Q32").Copy() | Out-Null
$ws.Range("R32").PasteSpecial(-4122) | Out-Null
$ws.Range("R32").Value = 27.9
$ws.Range("Q33").Copy() | Out-Null
$ws.Range("R33").PasteSpecial(-4122) | Out-Null
$ws.Range("R33").Value = 32.700000000000003
$ws.Range("Q34").Copy() | Out-Null
$ws.Range("R34").PasteSpecial(-4122) | Out-Null
$ws.Range("R34").Value = 22.8

$ws.Range("R3").Select() | Out-Null
$excel.CutCopyMode = 0

Write-Host "Edit complete"
